# "visualização de pedidos quase pronto"
# Adds a new request row (Maria Lima / Tecnologia) to the "pedidos_ao_rh"
# sheet, formats the new date cells, and leaves a stray formatted cell at
# N8 (mirroring the pre-existing stray cells at F5/I10), matching the
# author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 5: Maria Lima / Tecnologia, same date range as Pedro Pinheiro ---
$ws.Range("A5").Value = "Maria Lima"
$ws.Range("B5").Value = "Tecnologia"
$ws.Range("C5").Value = 44743
$ws.Range("D5").Value = 44772

# Give C5:D5 a plain date format (numFmtId 14, no alignment/font overrides).
# A helper cell + Copy/PasteSpecial(formats) is used so that both cells end
# up sharing a single newly-created style, instead of each Range.NumberFormat
# assignment minting its own separate style entry.
$ws.Range("Z1").NumberFormat = "mm-dd-yy"
$ws.Range("Z1").Copy()
$ws.Range("C5:D5").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("Z1").Clear()

# --- Stray formatted cell at N8 (same underline style as F5 / I10) ---
$ws.Range("N8").Font.Underline = $true

# Move the active selection to N8, like in the authored workbook.
$ws.Range("N8").Select()
